$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the sheet, pushing the existing data
# (and its formatting) down by one row.
$ws.Range("A1:D1").Insert()

# Fill in the new header row with the custom attribute names.
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "activity"
$ws.Range("C1").Value = "extraFieldA"
$ws.Range("D1").Value = "extraFieldB"

# Move the active selection to C4, matching the post-edit state.
$ws.Range("C4").Select()
